$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.399.51"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.773.16"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.37"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.32"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "3.771.34"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.485"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.86"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "4.409.69"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "3.778.21"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "69.528.32"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.55"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.26"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.44"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.45"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.21"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  +6.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.88"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.55"
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.54"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.97"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +6.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.339"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "464.98"
$ws.Range("E40").Value = "  +9.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.07"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.03"
$ws.Range("E42").Value = "  +9.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.81"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.42"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.58"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "2.952.35"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.35"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.97"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  -0.25%  "
